$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows identified by the diff (data rows removed from the Export sheet).
# Row numbers correspond to their position in the original sheet; deleting
# from bottom to top keeps the remaining row numbers stable during the loop.
$rowsToDelete = @(17, 13, 12, 9, 8, 7, 5)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
